$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Header row (row 1): add new columns F..M (bank, deposit_type, currency,
# owner, total, property_category, category, date, legislator_name,
# legislator_id, source_file, index). B1:D1 already hold bank/deposit_type/
# currency-style headers in the old sheet but need remapping; E1 becomes
# "owner", F1 becomes "total", G1..M1 are brand new.
$ws.Cells.Item(1,2).Value = "bank"
$ws.Cells.Item(1,3).Value = "deposit_type"
$ws.Cells.Item(1,4).Value = "currency"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "total"

$ws.Cells.Item(1,7).Value = "property_category"
$ws.Cells.Item(1,7).Font.Bold = $true
$ws.Cells.Item(1,7).Borders.LineStyle = 1
$ws.Cells.Item(1,7).HorizontalAlignment = -4108
$ws.Cells.Item(1,7).VerticalAlignment = -4160

$ws.Cells.Item(1,8).Value = "category"
$ws.Cells.Item(1,8).Font.Bold = $true
$ws.Cells.Item(1,8).Borders.LineStyle = 1
$ws.Cells.Item(1,8).HorizontalAlignment = -4108
$ws.Cells.Item(1,8).VerticalAlignment = -4160

$ws.Cells.Item(1,9).Value = "date"
$ws.Cells.Item(1,9).Font.Bold = $true
$ws.Cells.Item(1,9).Borders.LineStyle = 1
$ws.Cells.Item(1,9).HorizontalAlignment = -4108
$ws.Cells.Item(1,9).VerticalAlignment = -4160

$ws.Cells.Item(1,10).Value = "legislator_name"
$ws.Cells.Item(1,10).Font.Bold = $true
$ws.Cells.Item(1,10).Borders.LineStyle = 1
$ws.Cells.Item(1,10).HorizontalAlignment = -4108
$ws.Cells.Item(1,10).VerticalAlignment = -4160

$ws.Cells.Item(1,11).Value = "legislator_id"
$ws.Cells.Item(1,11).Font.Bold = $true
$ws.Cells.Item(1,11).Borders.LineStyle = 1
$ws.Cells.Item(1,11).HorizontalAlignment = -4108
$ws.Cells.Item(1,11).VerticalAlignment = -4160

$ws.Cells.Item(1,12).Value = "source_file"
$ws.Cells.Item(1,12).Font.Bold = $true
$ws.Cells.Item(1,12).Borders.LineStyle = 1
$ws.Cells.Item(1,12).HorizontalAlignment = -4108
$ws.Cells.Item(1,12).VerticalAlignment = -4160

$ws.Cells.Item(1,13).Value = "index"
$ws.Cells.Item(1,13).Font.Bold = $true
$ws.Cells.Item(1,13).Borders.LineStyle = 1
$ws.Cells.Item(1,13).HorizontalAlignment = -4108
$ws.Cells.Item(1,13).VerticalAlignment = -4160

# ---- Data rows (2..8): write column-by-column (B, C, D, E, F, G, H, I, J, K,
# L, M) so new shared strings are introduced in the same order the source
# workbook used.
$bank   = @{2="中華郵政股份有限公司中原大學郵局（中壢12支局）"; 3="臺灣銀行營業部"; 4="臺灣銀行群賢分行"; 5="臺灣銀行營業部"; 6="臺灣銀行營業部"; 7="中華郵政股份有限公司台北東門郵局(台北1支局）"; 8="臺灣銀行營業部"}
$dtype  = @{2="活期存款"; 3="活期存款"; 4="活期存款"; 5="活期存款"; 6="定期存款"; 7="活期存款"; 8="活期存款"}
$curr   = @{2="新臺幣"; 3="新臺幣"; 4="新臺幣"; 5="新臺幣"; 6="新臺幣"; 7="新臺幣"; 8="美金"}
$owner  = @{2="林治平"; 3="林治平"; 4="張曉風"; 5="張曉風"; 6="張曉風"; 7="張曉風"; 8="張曉風"}
$total  = @{2=327525; 3=44273; 4=359807; 5=3287844; 6=9100000; 7=1371654; 8=1619}
$pcat   = @{2="deposit"; 3="deposit"; 4="deposit"; 5="deposit"; 6="deposit"; 7="deposit"; 8="deposit"}
$cat    = @{2="normal"; 3="normal"; 4="normal"; 5="normal"; 6="normal"; 7="normal"; 8="normal"}
$date   = @{2="2013-03-17"; 3="2013-03-17"; 4="2013-03-17"; 5="2013-03-17"; 6="2013-03-17"; 7="2013-03-17"; 8="2013-03-17"}
$lname  = @{2="張曉風"; 3="張曉風"; 4="張曉風"; 5="張曉風"; 6="張曉風"; 7="張曉風"; 8="張曉風"}
$lid    = @{2=1748; 3=1748; 4=1748; 5=1748; 6=1748; 7=1748; 8=1748}
$sfile  = @{2="tmpd66d1"; 3="tmpd66d1"; 4="tmpd66d1"; 5="tmpd66d1"; 6="tmpd66d1"; 7="tmpd66d1"; 8="tmpd66d1"}
$index  = @{2=52; 3=53; 4=54; 5=55; 6=56; 7=57; 8=58}

for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,2).Value = $bank[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,3).Value = $dtype[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,4).Value = $curr[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,5).Value = $owner[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,6).Value = $total[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,7).Value = $pcat[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,8).Value = $cat[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,9).Value = $date[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,10).Value = $lname[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,11).Value = $lid[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,12).Value = $sfile[$r] }
for ($r = 2; $r -le 8; $r++) { $ws.Cells.Item($r,13).Value = $index[$r] }
